# Reverse the 0..15 / 15..0 numbering used throughout the "Hoja1" sheet of
# the A* map-way template. The header rows (1, 20, 41) listed the column
# index 15 -> 0 from left to right; they must now list 0 -> 15. The same
# applies to the "row index" columns A and T for each of the three 16-row
# blocks (rows 2-17, 21-36, 42-57). All the other cells on the sheet are
# formulas that depend on these values, so Excel will recompute them (and
# their cached <v> results) automatically once the literal values change.
#
# In addition, the view/selection stored with the sheet changes: the
# frozen/scrolled "topLeftCell" is cleared and the selected range moves
# from AG53:AI53 to A42:A57.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Build the 0..15 sequence once, as a 1x16 2D array (row vector) and as a
# 16x1 2D array (column vector), since Excel COM range assignment expects
# a 2D array matching the shape of the target range.
$rowVector = New-Object 'object[,]' 1,16
$colVector = New-Object 'object[,]' 16,1
for ($i = 0; $i -lt 16; $i++) {
    $rowVector[0, $i] = $i
    $colVector[$i, 0] = $i
}

# Header rows: columns B:Q and U:AJ each hold the 0..15 sequence.
foreach ($headerRow in 1, 20, 41) {
    $ws.Range("B$headerRow`:Q$headerRow").Value = $rowVector
    $ws.Range("U$headerRow`:AJ$headerRow").Value = $rowVector
}

# The three 16-row index blocks: column A and column T.
$blocks = @(
    @{ First = 2;  Last = 17 },
    @{ First = 21; Last = 36 },
    @{ First = 42; Last = 57 }
)

foreach ($block in $blocks) {
    $ws.Range("A$($block.First):A$($block.Last)").Value = $colVector
    $ws.Range("T$($block.First):T$($block.Last)").Value = $colVector
}

# Update the view: drop the old topLeftCell/selection (which scrolled to
# S34 and selected AG53:AI53) and select A42:A57 instead, from the top of
# the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A42:A57").Select()
